$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 39000
$ws.Range("J3").Value = 11700

# Row 4
$ws.Range("D4").Value = 7148277.97051587
$ws.Range("E4").Value = 7148278
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 23894783
$ws.Range("J4").Value = 7148278

# Row 5
$ws.Range("D5").Value = 3517249.673455395
$ws.Range("E5").Value = 3517250
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 11790829
$ws.Range("J5").Value = 3517250

# Row 6
$ws.Range("D6").Value = 10052394.79091849
$ws.Range("E6").Value = 10052395
$ws.Range("F6").Value = 2913077
$ws.Range("G6").Value = 873923.1
$ws.Range("H6").Value = 873923
$ws.Range("I6").Value = 57425760
$ws.Range("J6").Value = 10926318

# Row 7
$ws.Range("D7").Value = 3640201.440161694
$ws.Range("E7").Value = 3640201
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 12173302
$ws.Range("J7").Value = 3640201

# Row 9
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1792300
$ws.Range("J9").Value = 537690

# Row 10
$ws.Range("D10").Value = 3166283.696987225
$ws.Range("E10").Value = 3166284
$ws.Range("I10").Value = 10557199
$ws.Range("J10").Value = 3166284

# Row 11
$ws.Range("D11").Value = 4040336.403926175
$ws.Range("E11").Value = 4040336
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 17111679
$ws.Range("J11").Value = 4040336

# Row 12
$ws.Range("D12").Value = 9885758.627541639
$ws.Range("E12").Value = 9885759
$ws.Range("F12").Value = 1000000
$ws.Range("G12").Value = 300000
$ws.Range("H12").Value = 300000
$ws.Range("I12").Value = 50375382
$ws.Range("J12").Value = 10185759

# Row 13
$ws.Range("D13").Value = 2696374.240314188
$ws.Range("E13").Value = 2696374
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 9015215
$ws.Range("J13").Value = 2696374

# Row 14
$ws.Range("D14").Value = 8828358.856179321
$ws.Range("E14").Value = 8828359
$ws.Range("F14").Value = 1000000
$ws.Range("G14").Value = 300000
$ws.Range("H14").Value = 300000
$ws.Range("I14").Value = 33924822
$ws.Range("J14").Value = 9128359

# Row 17 (totals)
$ws.Range("D17").Value = 53524625.7
$ws.Range("E17").Value = 53524626
$ws.Range("F17").Value = 4913077
$ws.Range("G17").Value = 1473923.1
$ws.Range("H17").Value = 1473923
$ws.Range("I17").Value = 228100271
$ws.Range("J17").Value = 54998549
